$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove row 7 (CLR Investment Fund Public Ltd. (CSE:CLL)) entirely ---
$ws.Rows.Item(7).Delete()

# --- Row 2: Cyprus Asset Management (B2 code) ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '4'
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Value = -0.104
$ws.Range("G2").Value = 1.252442360296991
$ws.Range("H2").Value = 1.252442360296991
$ws.Range("I2").Value = 0.2999218444704963
$ws.Range("J2").Value = 0.2995139811146383
$ws.Range("K2").Value = 155.985
$ws.Range("L2").Value = 15.23886283704572
$ws.Range("M2").Value = 0.0
$ws.Range("N2").Value = 0.0
$ws.Range("O2").Value = 0.0
$ws.Range("P2").Value = 0.0
$ws.Range("Q2").Value = 0.0
$ws.Range("R2").Value = 0.0
$ws.Range("S2").Value = 0.0
$ws.Range("U2").Value = 9.996
$ws.Range("V2").Value = 0.0835716077251066
$ws.Range("W2").Value = -0.04026311188811188
$ws.Range("X2").Value = 0.05123888023031327
$ws.Range("Y2").Value = -0.09150199211842515
$ws.Range("Z2").Value = 0.05920185078079814
$ws.Range("AA2").Value = -0.5218194911185791
$ws.Range("AB2").Value = 0.05123888023031327
$ws.Range("AC2").Value = -0.5730583713488924
$ws.Range("AD2").Value = 10.7
$ws.Range("AE2").Value = 0.0
$ws.Range("AF2").Value = 10.7
$ws.Range("AG2").Value = 0.7039999999999988
$ws.Range("AH2").Value = 0.08211188703860026
$ws.Range("AI2").Value = 0.03037184217996026
$ws.Range("AJ2").Value = 0.005851355619462397
$ws.Range("AK2").Value = 0.002056651397588106
$ws.Range("AL2").Value = 0.427
$ws.Range("AM2").Value = -0.242
$ws.Range("AN2").Value = 2.196674194210634
$ws.Range("AO2").Value = 7.189695550351288
$ws.Range("AP2").Value = 0.1445288441798396
$ws.Range("AQ2").Value = -12.68595041322314
$ws.Range("E2").ClearContents()
$ws.Range("T2").ClearContents()

# --- Row 3: Demetra Holdings Plc (CSE:DEM) ---
$ws.Range("B3").Value = 'Demetra Holdings Plc (CSE:DEM)'
$ws.Range("G3").Value = 1.673101673101673
$ws.Range("H3").Value = 1.673101673101673
$ws.Range("I3").Value = 0.7284427284427285
$ws.Range("J3").Value = 0.7244802948787769
$ws.Range("K3").Value = 157.3
$ws.Range("L3").Value = 20.24453024453025
$ws.Range("U3").Value = 1.66
$ws.Range("V3").Value = 0.01655034895314058
$ws.Range("W3").Value = 1.104634831460674
$ws.Range("X3").Value = 0.05515366354194147
$ws.Range("Y3").Value = 1.049481167918733
$ws.Range("Z3").Value = 0.05154912757911496
$ws.Range("AA3").Value = 0.0373463271492609
$ws.Range("AB3").Value = 0.0536957920921542
$ws.Range("AC3").Value = -0.0163494649428933
$ws.Range("AD3").Value = 10.7
$ws.Range("AF3").Value = 10.7
$ws.Range("AG3").Value = 9.04
$ws.Range("AH3").Value = 0.09639639639639638
$ws.Range("AI3").Value = 0.03342705404561075
$ws.Range("AJ3").Value = 0.0826778854947869
$ws.Range("AK3").Value = 0.02838839341791232
$ws.Range("AL3").Value = 0.426
$ws.Range("AM3").Value = 0.426
$ws.Range("AN3").Value = 1.854419410745234
$ws.Range("AO3").Value = 13.28638497652582
$ws.Range("AP3").Value = 1.566724436741768
$ws.Range("AQ3").Value = 13.28638497652582
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# --- Row 4: Argo Group Limited (AIM:ARGO) ---
$ws.Range("D4").Value = -0.104
$ws.Range("G4").Value = -0.046875
$ws.Range("H4").Value = -0.046875
$ws.Range("I4").Value = -0.23671875
$ws.Range("J4").Value = -0.23671875
$ws.Range("K4").Value = -0.276
$ws.Range("L4").Value = -0.07187500000000001
$ws.Range("M4").Value = -0.0
$ws.Range("N4").Value = -0.0
$ws.Range("O4").Value = 0.0
$ws.Range("R4").Value = 0.0
$ws.Range("S4").Value = 0.0
$ws.Range("U4").Value = 1.6
$ws.Range("V4").Value = 0.1626016260162602
$ws.Range("W4").Value = -0.01254545454545455
$ws.Range("X4").Value = 0.05123888023031327
$ws.Range("Y4").Value = -0.06378433477576781
$ws.Range("Z4").Value = 0.1843494959193471
$ws.Range("AA4").Value = -0.04363898223715795
$ws.Range("AB4").Value = 0.05123888023031327
$ws.Range("AC4").Value = -0.09487786246747122
$ws.Range("AE4").Value = 0.0
$ws.Range("AF4").Value = 0.0
$ws.Range("AG4").Value = -1.6
$ws.Range("AH4").Value = 0.0
$ws.Range("AI4").Value = 0.0
$ws.Range("AJ4").Value = -0.1941747572815534
$ws.Range("AK4").Value = -0.08080808080808083
$ws.Range("AM4").Value = -0.493
$ws.Range("AP4").Value = 1.779755283648498
$ws.Range("AQ4").Value = 1.843813387423935
$ws.Range("T4").ClearContents()

# --- Row 5: Interfund Investments Plc. (CSE:INF) ---
$ws.Range("I5").Value = 1.229357798165138
$ws.Range("J5").Value = 1.229357798165138
$ws.Range("K5").Value = -0.707
$ws.Range("L5").Value = 0.6486238532110091
$ws.Range("U5").Value = 6.61
$ws.Range("V5").Value = 0.8529032258064516
$ws.Range("W5").Value = -0.06798076923076922
$ws.Range("X5").Value = 0.05123888023031327
$ws.Range("Y5").Value = -0.1192196494610825
$ws.Range("Z5").Value = -0.8134328358208956
$ws.Range("AA5").Value = -1.0
$ws.Range("AB5").Value = 0.05123888023031327
$ws.Range("AC5").Value = -1.051238880230313
$ws.Range("AG5").Value = -6.61
$ws.Range("AJ5").Value = -5.79824561403509
$ws.Range("AK5").Value = -2.248299319727891
$ws.Range("AL5").Value = 0.0
$ws.Range("AM5").Value = -0.176
$ws.Range("AQ5").Value = 7.613636363636364
$ws.Range("AO5").ClearContents()

# --- Row 6: Actibond Growth Fund Public Co. Ltd. (CSE:ACT) ---
$ws.Range("I6").Value = 1.200704225352113
$ws.Range("J6").Value = 1.200704225352113
$ws.Range("K6").Value = -0.332
$ws.Range("L6").Value = 1.169014084507042
$ws.Range("U6").Value = 0.126
$ws.Range("V6").Value = 0.07325581395348837
$ws.Range("W6").Value = -0.2075
$ws.Range("X6").Value = 0.05123888023031327
$ws.Range("Y6").Value = -0.2587388802303133
$ws.Range("AB6").Value = 0.05123888023031327
$ws.Range("AG6").Value = -0.126
$ws.Range("AJ6").Value = -0.07904642409033878
$ws.Range("AK6").Value = -0.1120996441281139
$ws.Range("AL6").Value = 0.001
$ws.Range("AM6").Value = 0.001
$ws.Range("AO6").Value = -341.0
$ws.Range("AQ6").Value = -341.0
$ws.Range("Z6").ClearContents()
$ws.Range("AA6").ClearContents()
$ws.Range("AC6").ClearContents()

